$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case TC009 "Verify Edit Internship Title" (modeled on TC008,
# "Verify Edit Job Title"), appended as rows 73-81.
#
# Columns: testCaseId, testCaseName, id, action, selector, data
$rows = @(
    @("TC009", "Verify Edit Internship Title", "1", "login_to_employer_portal", "", ""),
    @("TC009", "Verify Edit Internship Title", "2", "click", "RECENT_INTERNSHIPS_LINK", ""),
    @("TC009", "Verify Edit Internship Title", "3", "waitfor", "EDIT_JOB_BTN", ""),
    @("TC009", "Verify Edit Internship Title", "4", "click", "EDIT_JOB_BTN", ""),
    @("TC009", "Verify Edit Internship Title", "5", "waitfor", "INTERNSHIP_TITLE_INPUT", ""),
    @("TC009", "Verify Edit Internship Title", "6", "type", "INTERNSHIP_TITLE_INPUT", "demo internship {{TIMESTAMP}} edit"),
    @("TC009", "Verify Edit Internship Title", "7", "click", "UPDATE_INTERNSHIP_BTN", ""),
    @("TC009", "Verify Edit Internship Title", "8", "click", "ALERT_OK_BTN", ""),
    @("TC009", "Verify Edit Internship Title", "9", "verify_text", "css:body", "demo internship {{TIMESTAMP}} edit")
)

$startRow = 73
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    for ($c = 1; $c -le 6; $c++) {
        $val = $rows[$i][$c - 1]
        # Every column in this sheet is stored as text, even purely numeric
        # ids ("id" column) and blanks. Plain text assignments stay text on
        # their own, but a numeric-looking or empty value needs a leading
        # apostrophe so it is kept as text instead of being coerced to a
        # number/blank - so only use that trick where it's actually needed.
        if ($val -match '^-?[0-9]+(\.[0-9]+)?$' -or $val -eq "") {
            $ws.Cells.Item($r, $c).Value = "'" + $val
        } else {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

# The sheet view no longer pins an explicit right-to-left flag.
$ws.DisplayRightToLeft = $false
